$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new pair of rows (Primera/Segunda for a new date, 2022-02-14)
# right before the current row 221, pushing the existing data down by two rows.
$ws.Range("A221:A222").EntireRow.Insert()

# New row 221: Apio, Americana (o), Primera, fecha 2022-02-14 (44606)
$ws.Range("A221").Value = 8
$ws.Range("B221").Value = "Terminal La Palmera de La Serena"
$ws.Range("C221").Value = "Coquimbo"
$ws.Range("D221").Value = 44606
$ws.Range("E221").Value = 4
$ws.Range("F221").Value = 100112017
$ws.Range("G221").Value = "Apio"
$ws.Range("H221").Value = "Americana (o)"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 2000
$ws.Range("K221").Value = 8000
$ws.Range("L221").Value = 9000
$ws.Range("M221").Value = 8500
$ws.Range("N221").Value = "`$/docena de matas"
$ws.Range("O221").Value = "Provincia del Elquí"
$ws.Range("P221").Value = 1417
$ws.Range("Q221").Value = 6
$ws.Range("R221").Value = "Hortaliza"

# New row 222: Apio, Americana (o), Segunda, fecha 2022-02-14 (44606)
$ws.Range("A222").Value = 8
$ws.Range("B222").Value = "Terminal La Palmera de La Serena"
$ws.Range("C222").Value = "Coquimbo"
$ws.Range("D222").Value = 44606
$ws.Range("E222").Value = 4
$ws.Range("F222").Value = 100112017
$ws.Range("G222").Value = "Apio"
$ws.Range("H222").Value = "Americana (o)"
$ws.Range("I222").Value = "Segunda"
$ws.Range("J222").Value = 1360
$ws.Range("K222").Value = 6000
$ws.Range("L222").Value = 7000
$ws.Range("M222").Value = 6500
$ws.Range("N222").Value = "`$/docena de matas"
$ws.Range("O222").Value = "Provincia del Elquí"
$ws.Range("P222").Value = 1083
$ws.Range("Q222").Value = 6
$ws.Range("R222").Value = "Hortaliza"
